$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "06/08/2025"
$ws.Range("A23").Style = "Normal"
$ws.Range("B23").Value = "A. Lima"
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = "Sporting Cristal"
$ws.Range("F23").Value = "D"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1.24
$ws.Range("L23").Value = 0.25
$ws.Range("M23").Value = 17
$ws.Range("N23").Value = 5
$ws.Range("O23").Value = 1
$ws.Range("P23").Value = 2
